$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cell = $ws.Range("A9")
$cell | Get-Member
